# Update cryptocurrency price/volume data in the worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.477.86"
$ws.Range("E2").Value = "  +0.12%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.805.92"
$ws.Range("E3").Value = "  +0.16%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "224.80"
$ws.Range("E5").Value = "  -1.20%  "
$ws.Range("E6").Value = "  +4.66%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "38.92"
$ws.Range("E8").Value = "  +7.07%  "
$ws.Range("E9").Value = "  -3.11%  "
$ws.Range("E10").Value = "  -3.12%  "
$ws.Range("E11").Value = "  +1.96%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.067.60"
$ws.Range("E12").Value = "  +0.20%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.07"
$ws.Range("E13").Value = "  -4.13%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.805.81"
$ws.Range("E14").Value = "  -0.21%  "
$ws.Range("E15").Value = "  -2.02%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "34.452.88"
$ws.Range("E16").Value = "  +0.12%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.37"
$ws.Range("E17").Value = "  -2.25%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.14"
$ws.Range("E18").Value = "  -2.74%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "240.94"
$ws.Range("E19").Value = "  -1.64%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0768"
$ws.Range("E20").Value = "  -2.75%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.19"
$ws.Range("E21").Value = "  -3.36%  "
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.09"
$ws.Range("E23").Value = "  -1.88%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.18"
$ws.Range("E24").Value = "  +0.69%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "171.01"
$ws.Range("E25").Value = "  -0.93%  "
$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "17.61"
$ws.Range("E26").Value = "  +4.18%  "
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.68"
$ws.Range("E27").Value = "  -3.62%  "
$ws.Range("E28").Value = "  +3.63%  "
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.22"
$ws.Range("E30").Value = "  -1.29%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.76"
$ws.Range("E31").Value = "  -2.06%  "
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0514"
$ws.Range("E32").Value = "  -2.67%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.85"
$ws.Range("E33").Value = "  -3.90%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.82"
$ws.Range("E34").Value = "  +0.87%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.641"
$ws.Range("E35").Value = "  -4.39%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.307.18"
$ws.Range("E36").Value = "  -6.37%  "
$ws.Range("E37").Value = "  -0.81%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.32"
$ws.Range("E38").Value = "  -4.43%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0186"
$ws.Range("E39").Value = "  -1.93%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "82.34"
$ws.Range("E40").Value = "  -0.11%  "
$ws.Range("E41").Value = "  +0.64%  "
$ws.Range("E42").Value = "  -0.28%  "
$ws.Range("E43").Value = "  +1.52%  "
$ws.Range("E44").Value = "  -1.15%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.94"
$ws.Range("E45").Value = "  +4.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0512"
$ws.Range("E46").Value = "  +1.94%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.967.94"
$ws.Range("E47").Value = "  +0.16%  "
$ws.Range("E48").Value = "  -4.44%  "
$ws.Range("E49").Value = "  -0.15%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "102.63"
$ws.Range("E50").Value = "  -1.38%  "
$ws.Range("E51").Value = "  -6.39%  "
